$wb = $excel.ActiveWorkbook
$wsProperties = $wb.Worksheets.Item(1)
$wsLeases = $wb.Worksheets.Item(2)

# --- Leases sheet: replace tenant record (Kuldeep Yadav -> Virat Kohli) ---
# Write Last Name before First Name so the shared-string table picks up the
# same insertion order as the target workbook (Kohli, then Virat).
$wsLeases.Range("B2").Value = "Kohli"
$wsLeases.Range("A2").Value = "Virat"
$wsLeases.Range("C2").Value = "virat.kohli@nomail.com"

# Swap out the mailto hyperlink so it points at the new email address while
# keeping the existing "Hyperlink" cell styling.
$wsLeases.Hyperlinks.Delete()
$linkCell = $wsLeases.Range("C2")
$linkCell.Hyperlinks.Add($linkCell, "mailto:virat.kohli@nomail.com")
$linkCell.Style = "Hyperlink"

# --- View state: Leases becomes the active/selected tab ---
$null = $wsProperties.Activate()
$null = $wsProperties.Range("D16").Select()

$null = $wsLeases.Activate()
$null = $wsLeases.Range("C4").Select()
